# No more archive files in ingests
# Removes the three "archive" (zip) related records from the
# "Records to Create" sheet: "Broken zip", "Broken zip file" and
# "DOOM shareware version". Deleting the entire rows shifts the
# remaining rows up and updates dependent ranges (data validation,
# shared strings, etc.) automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Records to Create")

# Row 12: "Broken zip"
# Row 18: "Broken zip file"
# Row 21: "DOOM shareware version"
# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows.Item(21).Delete()
$ws.Rows.Item(18).Delete()
$ws.Rows.Item(12).Delete()

# Update the active selection on the sheet to match the saved view.
$ws.Activate()
$ws.Range("R19").Select()
